$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 3, shifting existing rows 3-10 down to 5-12.
$ws.Range("A3:A4").EntireRow.Insert()

# New row 3: Terminal Hortofrutícola Agro Chillán, Fecha 2023-12-21, Primera
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 45281
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101004
$ws.Cells.Item(3, 10).Value = "Frambuesa"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 120
$ws.Cells.Item(3, 14).Value = 7500
$ws.Cells.Item(3, 15).Value = 8000
$ws.Cells.Item(3, 16).Value = 7750
$ws.Cells.Item(3, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value = "Región de Ñuble"
$ws.Cells.Item(3, 19).Value = 3875
$ws.Cells.Item(3, 20).Value = 2

# New row 4: Terminal Hortofrutícola Agro Chillán, Fecha 2023-12-21, Segunda
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 45281
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100101
$ws.Cells.Item(4, 8).Value = "Berries"
$ws.Cells.Item(4, 9).Value = 100101004
$ws.Cells.Item(4, 10).Value = "Frambuesa"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 120
$ws.Cells.Item(4, 14).Value = 6000
$ws.Cells.Item(4, 15).Value = 6500
$ws.Cells.Item(4, 16).Value = 6250
$ws.Cells.Item(4, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(4, 18).Value = "Región de Ñuble"
$ws.Cells.Item(4, 19).Value = 3125
$ws.Cells.Item(4, 20).Value = 2
